# Actualizar cabanas y calendario
# Update the "Hoja1" (users) report sheet with refreshed data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Data rows (ID, Nombre, Apellido Paterno, Apellido Materno, Alias, Email, Tipo Usuario)
$data = @(
    @(14, "Prueba",          "Prueba",    "Prueba",      "Prueba",         "aaaa@gmail.com",               "Administrador"),
    @(11, "Yomer",            "asies",     "asies",       "asies",          "asies@gmail.com",               "Administrador"),
    @(10, "Ramiro",           "De Jesus",  "Hernandez",   "RamboBernabe",   "ramironchis@gmail.com",         "Administrador"),
    @(9,  "Gilberta",         "Olivares",  "Cruz",        "ElGil",          "gilimemo@gmail.com",            "Administrador"),
    @(3,  "Erick Jonathan ",  "Bautista",  "Perez",       "f4k3r",          "erick.bautista57@hotmail.com",  "Administrador"),
    @(6,  "Ramiro de Jesús",  "Hernández", "Bernabé ",    "nalgoncito",     "ramironalgon@gmail.com",        "Administrador"),
    @(7,  "José Manuel",      "Jimeno ",   "Islas ",      "ChemigodElite ", "elitefortgod@gmail.com",        "Administrador"),
    @(4,  "Gilberto",         "Cruz",      "Olivares",    "C0Gil",          "gil123@gmail.com",              "Administrador"),
    @(17, "Luis2",            "Lopez2",    "Delgado2",    "Cerre2",         "luis2@gmail.com",               "Supervisor"),
    @(16, "Luis",             "Lopez",     "Delgado",     "Cerre",          "luis@gmail.com",                "Supervisor"),
    @(15, "Teofilito",        "Hernandez", "Bernabe",     "Teo",            "teo@gmail.com",                 "Supervisor"),
    @(13, "Misterbist",       "Señor",     "Bestia",      "MrB",            "pokimike@gmail.com",            "Supervisor"),
    @(8,  "Luis Eduardo",     "Bautista",  "Perez",       "elsanto",        "luisbautista@example.com",      "Supervisor"),
    @(5,  "Juan Miguel",      "Sanchez",   "Aguilar",     "pokimike",       "pokimike@gmail.com",            "Supervisor")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
